# The authored change swaps the presentation's theme color palette from the
# "Integral" (Red Violet) scheme to the default "Office" scheme.
#
# Index order for ThemeColorScheme.Colors(n) (verified against the deck's
# existing Integral palette): 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2
# 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink.
#
# PowerPoint COM packs RGB() as r + g*256 + b*65536 (Windows COLORREF / BGR
# byte order once serialized), so we build that value from the familiar hex
# triplet for each OOXML srgbClr.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeRGB {
    param($Index, $HexColor)
    $r = [Convert]::ToInt32($HexColor.Substring(0,2), 16)
    $g = [Convert]::ToInt32($HexColor.Substring(2,2), 16)
    $b = [Convert]::ToInt32($HexColor.Substring(4,2), 16)
    $tcs.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeRGB 1  "000000"   # dk1
Set-ThemeRGB 2  "FFFFFF"   # lt1
Set-ThemeRGB 3  "44546A"   # dk2
Set-ThemeRGB 4  "E7E6E6"   # lt2
Set-ThemeRGB 5  "5B9BD5"   # accent1
Set-ThemeRGB 6  "ED7D31"   # accent2
Set-ThemeRGB 7  "A5A5A5"   # accent3
Set-ThemeRGB 8  "FFC000"   # accent4
Set-ThemeRGB 9  "4472C4"   # accent5
Set-ThemeRGB 10 "70AD47"   # accent6
Set-ThemeRGB 11 "0563C1"   # hlink
Set-ThemeRGB 12 "954F72"   # folHlink
